$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.645.27"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.032.29"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "382.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "3.509.48"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "3.036.05"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -11.59%  "
$ws.Range("D19").Value = "51.662.05"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.48%  "
$ws.Range("E27").Value = "  +6.09%  "
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0447"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.294"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("D49").Value = "2.028.57"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").Value = "3.332.96"
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.516"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.05%  "
